# Towers vs. Enemies chart.xlsx -- "Enemies shoot now, but it's in progress"
#
# The enemy stat-block (rows 34-45) is being reworked to look like the tower
# stat-block: "Speed" -> "Reload", "Damage(per sec)" -> "Damage", and the
# numbers in the Reload/Damage columns are bumped up (to real per-second
# scale) instead of the old small placeholder integers. A few tower rows
# also get tweaked (cost/damage numbers, and the "slow" Fire Rate for
# Police/Fire Department becomes "very slow"). Finally the header label
# "Damage (per sec)" on the tower table becomes just "Damage".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Tower table header row (row 23) ----
$ws.Range("C23").Value = "Damage"

# ---- Tower stat tweaks (rows 24-32) ----
$ws.Range("F24").Value = 2
$ws.Range("C27").Value = 8
$ws.Range("C28").Value = 6
$ws.Range("C29").Value = 1000
$ws.Range("H29").Value = "very slow"
$ws.Range("C31").Value = 25

# ---- Enemy table header row (row 34) ----
$ws.Range("B34").Value = "Reload"
$ws.Range("C34").Value = "Damage"

# ---- Enemy stat tweaks (rows 35-45) ----
$ws.Range("B35").Value = 25
$ws.Range("C35").Value = 2

$ws.Range("B36").Value = 100
$ws.Range("C36").Value = 10

$ws.Range("B37").Value = 500
$ws.Range("C37").Value = 50

$ws.Range("B38").Value = 150
$ws.Range("C38").Value = 15

$ws.Range("B39").Value = 50

$ws.Range("B40").Value = 25

$ws.Range("B41").Value = 50
$ws.Range("C41").Value = 5

$ws.Range("B42").Value = 25
$ws.Range("C42").Value = 5

$ws.Range("B43").Value = 5

$ws.Range("B44").Value = 50
$ws.Range("C44").Value = 5

$ws.Range("B45").Value = 25

# ---- View state: scrolled down a bit, selection left on B45 ----
$ws.Activate()
try {
    $excel.ActiveWindow.ScrollRow = 13
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
$ws.Range("B45").Select()
